$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# LeBron's second Cleveland stint (rows 844-1144, i.e. game numbers 843-1143)
# needs to be relabeled "CLET" to distinguish it from his first Cleveland stint ("CLE").
for ($r = 844; $r -le 1144; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq "CLE") {
        $cell.Value = "CLET"
    }
}

# Update the saved view position/selection to match where the user was working.
$ws.Range("D1144").Select()
